# Apply edits to "Perc Decline in Battery Cost per Doubling of Cap.xlsx"
# Source is switched from Bloomberg New Energy Finance (2018) to an MIT paper (2021),
# the illustrative picture on the About sheet is removed, and the PDiBCpDoC value
# becomes a computed average instead of a hard-coded constant.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "About"
$ws2 = $wb.Worksheets.Item(2)   # "PDiBCpDoC"

# ---------------------------------------------------------------------------
# About sheet
# ---------------------------------------------------------------------------

# Drop the "last updated" date stamp that used to sit in C1
$ws1.Range("C1").Clear()

# Row 6: source URL
$ws1.Range("B6").Value = "https://pubs.rsc.org/en/content/articlepdf/2021/ee/d0ee02681f?page=search"

# Row 5: source title
$ws1.Range("B5").Value = "Re-examining rates of lithium-ion battery technology improvement and cost decline"

# Row 3: source name
$ws1.Range("B3").Value = "Massachusetts Institute of Technology"

# Row 7: source location note
$ws1.Range("B7").Value = "Abstract"

# Row 9 (new): methodology note
$ws1.Range("A9").Value = "Note: We take the average of learning rates quoted in the Abstract (20%-27%)"

# Row 4: source year
$ws1.Range("B4").Value = 2021

# Row 8: the old "note" text is gone, but keep the (italic) cell, now empty
$ws1.Range("C8").Value = ""

# Remove the embedded chart picture that illustrated the old BNEF source
while ($ws1.Shapes.Count -gt 0) {
    $ws1.Shapes.Item(1).Delete()
}

# ---------------------------------------------------------------------------
# PDiBCpDoC sheet
# ---------------------------------------------------------------------------

# Relabel the value column header
$ws2.Range("B1").Value = "Perc Decline per Doubling (dimensionless)"

# Replace the hard-coded 0.18 with a computed average of the two quoted rates
$ws2.Range("B2").Formula = "=AVERAGE(0.2,0.27)"

# ---------------------------------------------------------------------------
# Restore on-screen selections to match the saved view state
# ---------------------------------------------------------------------------
$ws2.Range("I4").Select()
$ws1.Range("A10").Select()
